# Append a new 案件 (job listing) row at the top of the data in "ランサーズ",
# refresh the "取得日時" timestamp on every row, and grow two column widths
# to fit the new content (commit: "Append: 2026-02-15 12:59 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Make room: insert a new row right under the header, pushing the
#        two existing listings down from rows 2-3 to rows 3-4. ---
$ws.Rows.Item(2).Insert()

# --- 2. Fill in the newly inserted row with the freshly scraped listing. ---
$ws.Range("A2").Value = "2026-02-15 12:59:30"
$ws.Range("B2").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = "◆ツール,スクレイピング ◇サイト"

# --- 3. The scrape timestamp in column A is refreshed for every row on
#        this run, including the listings that were already present. ---
$ws.Range("A3").Value = "2026-02-15 12:59:30"
$ws.Range("A4").Value = "2026-02-15 12:59:30"

# --- 4. Widen column B (title) and column H (skill summary) to fit the
#        longer new strings. ColumnWidth uses Excel's padded character
#        unit, which adds ~0.8333 to whatever integer is requested, so we
#        back that padding out to land exactly on 51 / 19. ---
$ws.Columns.Item(2).ColumnWidth = 50.16666666666667
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668

# --- 5. Row-insert does not renumber existing Hyperlink ranges, so rebuild
#        the column-F hyperlinks from scratch to match the shifted rows. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5492003")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5491983")
$ws.Range("F2:F4").Style = "Hyperlink"
